$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.018720333333333
$ws.Range("H2").Value = 6.056161
$ws.Range("I2").Value = 0.02546714762636305
$ws.Range("J2").Value = 0.02546714762636305
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 3.873735232862111
$ws.Range("R2").Value = 34.863617095759
$ws.Range("S2").Value = 0.0001661481591288336
$ws.Range("T2").Value = 0.0001661481591288336
$ws.Range("G3").Value = 2.018720333333333
$ws.Range("H3").Value = 6.056161
$ws.Range("I3").Value = 0.02546714762636305
$ws.Range("J3").Value = 0.02546714762636305
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 365.9705617818905
$ws.Range("R3").Value = 3293.735056037014
$ws.Range("S3").Value = 0.01569682269959378
$ws.Range("T3").Value = 0.01569682269959378
$ws.Range("G4").Value = 2.018720333333333
$ws.Range("H4").Value = 6.056161
$ws.Range("I4").Value = 0.02546714762636305
$ws.Range("J4").Value = 0.02546714762636305
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 58.63912677530412
$ws.Range("R4").Value = 527.752140977737
$ws.Range("S4").Value = 0.002515087475258503
$ws.Range("T4").Value = 0.002515087475258504
$ws.Range("G5").Value = 2.018720333333333
$ws.Range("H5").Value = 6.056161
$ws.Range("I5").Value = 0.02546714762636305
$ws.Range("J5").Value = 0.02546714762636305
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 165.2817286980085
$ws.Range("R5").Value = 1487.535558282076
$ws.Range("S5").Value = 0.00708908929238193
$ws.Range("T5").Value = 0.00708908929238193
$ws.Range("I6").Value = 0.5652548410284588
$ws.Range("J6").Value = 0.5652548410284589
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 85.9793026436591
$ws.Range("R6").Value = 773.8137237929318
$ws.Range("S6").Value = 0.003687733414570544
$ws.Range("T6").Value = 0.003687733414570545
$ws.Range("I7").Value = 0.5652548410284588
$ws.Range("J7").Value = 0.5652548410284589
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.3483980675765177
$ws.Range("T7").Value = 0.3483980675765178
$ws.Range("I8").Value = 0.5652548410284588
$ws.Range("J8").Value = 0.5652548410284589
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 1301.521896747875
$ws.Range("R8").Value = 11713.69707073087
$ws.Range("S8").Value = 0.05582350217847858
$ws.Range("T8").Value = 0.05582350217847859
$ws.Range("I9").Value = 0.5652548410284588
$ws.Range("J9").Value = 0.5652548410284589
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 3668.502599929516
$ws.Range("R9").Value = 33016.52339936564
$ws.Range("S9").Value = 0.157345537858892
$ws.Range("T9").Value = 0.157345537858892
$ws.Range("G10").Value = 31.416885
$ws.Range("H10").Value = 94.25065499999999
$ws.Range("I10").Value = 0.3963394210897649
$ws.Range("J10").Value = 0.396339421089765
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 60.28605960010499
$ws.Range("R10").Value = 542.5745364009449
$ws.Range("S10").Value = 0.00258572597804728
$ws.Range("T10").Value = 0.002585725978047281
$ws.Range("G11").Value = 31.416885
$ws.Range("H11").Value = 94.25065499999999
$ws.Range("I11").Value = 0.3963394210897649
$ws.Range("J11").Value = 0.396339421089765
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 5695.51654235433
$ws.Range("R11").Value = 51259.64888118897
$ws.Range("S11").Value = 0.2442860784010831
$ws.Range("T11").Value = 0.2442860784010831
$ws.Range("G12").Value = 31.416885
$ws.Range("H12").Value = 94.25065499999999
$ws.Range("I12").Value = 0.3963394210897649
$ws.Range("J12").Value = 0.396339421089765
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 912.587381214015
$ws.Range("R12").Value = 8213.286430926135
$ws.Range("S12").Value = 0.03914173383524813
$ws.Range("T12").Value = 0.03914173383524815
$ws.Range("G13").Value = 31.416885
$ws.Range("H13").Value = 94.25065499999999
$ws.Range("I13").Value = 0.3963394210897649
$ws.Range("J13").Value = 0.396339421089765
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 2572.24191848922
$ws.Range("R13").Value = 23150.17726640298
$ws.Range("S13").Value = 0.1103258828753865
$ws.Range("T13").Value = 0.1103258828753865
$ws.Range("G14").Value = 1.025611333333333
$ws.Range("H14").Value = 3.076834
$ws.Range("I14").Value = 0.01293859025541314
$ws.Range("J14").Value = 0.01293859025541314
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 1.968052083071777
$ws.Range("R14").Value = 17.712468747646
$ws.Range("S14").Value = 0.00008441161076216526
$ws.Range("T14").Value = 0.00008441161076216528
$ws.Range("G15").Value = 1.025611333333333
$ws.Range("H15").Value = 3.076834
$ws.Range("I15").Value = 0.01293859025541314
$ws.Range("J15").Value = 0.01293859025541314
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 185.9314287532351
$ws.Range("R15").Value = 1673.382858779116
$ws.Range("S15").Value = 0.007974774411393937
$ws.Range("T15").Value = 0.007974774411393937
$ws.Range("G16").Value = 1.025611333333333
$ws.Range("H16").Value = 3.076834
$ws.Range("I16").Value = 0.01293859025541314
$ws.Range("J16").Value = 0.01293859025541314
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 29.79162195201977
$ws.Range("R16").Value = 268.124597568178
$ws.Range("S16").Value = 0.001277790774857128
$ws.Range("T16").Value = 0.001277790774857128
$ws.Range("G17").Value = 1.025611333333333
$ws.Range("H17").Value = 3.076834
$ws.Range("I17").Value = 0.01293859025541314
$ws.Range("J17").Value = 0.01293859025541314
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 83.9714205809271
$ws.Range("R17").Value = 755.7427852283439
$ws.Range("S17").Value = 0.00360161345839991
$ws.Range("T17").Value = 0.003601613458399911
